# "Health + heroism added" -- the spellbook list gets re-sorted alphabetically
# by Discipline (column A), an AutoFilter is turned on for the table, and the
# current selection / scroll position are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A1:G15")
$sortKey   = $ws.Range("A1")

# Sort rows 2-15 alphabetically (ascending) by the Discipline column, keeping
# the header row (row 1) fixed.
$dataRange.Sort($sortKey, 1, $null, $null, 1, $null, 1, 1, 2) | Out-Null

# The row heights in the source file are "auto" heights driven by how much
# text each discipline's longest title needs to wrap into -- they travel
# with the data when Excel re-sorts the rows, so re-apply them explicitly
# for the rows whose content -- and therefore needed height -- changed.
$ws.Range("A2:G2").RowHeight   = 35.05
$ws.Range("A3:G3").RowHeight   = 23.85
$ws.Range("A4:G4").RowHeight   = 23.85
$ws.Range("A5:G5").RowHeight   = 23.85
$ws.Range("A6:G6").RowHeight   = 23.85
$ws.Range("A7:G7").RowHeight   = 23.85
$ws.Range("A8:G8").RowHeight   = 23.85
$ws.Range("A9:G9").RowHeight   = 23.85
$ws.Range("A10:G10").RowHeight = 13.8
$ws.Range("A11:G11").RowHeight = 23.85
$ws.Range("A12:G12").RowHeight = 23.85
$ws.Range("A13:G13").RowHeight = 23.85
$ws.Range("A14:G14").RowHeight = 23.85
$ws.Range("A15:G15").RowHeight = 23.85

# Turn the AutoFilter on for the table.
$dataRange.AutoFilter() | Out-Null

# Excel records the AutoFilter's range as a hidden, sheet-scoped defined
# name, "_FilterDatabase".
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$G`$15")
$filterName.Visible = $false

# Update the view: scrolled so column C is left-most, and K12 selected.
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
$ws.Range("K12").Select() | Out-Null
